$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.759.33"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.02%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.139.98"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -7.98%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "565.23"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.96%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.05"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -4.11%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.615"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.96%  "

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.04%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.137.19"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -7.89%  "

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -5.58%  "

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -5.87%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.392"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -4.42%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.685.55"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -7.89%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.135"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.96%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.89"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -7.36%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.674.23"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.16%  "

# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -5.39%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.140.23"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -8.06%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.66"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.70%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.78"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -6.77%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "354.23"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.18%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.23"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -4.07%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.27%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.44"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -6.23%  "

# Row 25
$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.497"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -6.59%  "

# Row 26
$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000117"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -6.87%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.51"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.35%  "

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.12%  "

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.06%  "

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.10%  "

# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.60%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.34"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -6.70%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.84"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -6.13%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.61"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.85%  "

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -4.63%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.43"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -6.49%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "153.60"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -5.65%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.829"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.26%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "25.98"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -5.43%  "

# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.58%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.50"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.36%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.662.94"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.01%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.15"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -5.96%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.01"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -4.19%  "

# Row 45
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0655"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -4.33%  "

# Row 46
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.13"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.81%  "

# Row 47
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.15"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.28%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "323.18"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.34%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0271"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.48%  "

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.14%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.999"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.06%  "
